# Auto-generated edit script applying the cryptos.xlsx diff
# Updates Price (D) and Volume(1h) (E) columns, and swaps the
# ThetaToken / InjectiveProtocol rows (48 and 50).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "63.116.29"
Set-TextValue "E2" "  -2.10%  "
Set-TextValue "D3" "3.127.86"
Set-TextValue "E3" "  -0.58%  "
Set-TextValue "E4" "  -0.04%  "
Set-TextValue "D5" "594.22"
Set-TextValue "E5" "  -2.77%  "
Set-TextValue "D6" "136.39"
Set-TextValue "E6" "  -5.13%  "
Set-TextValue "E7" "  +0.01%  "
Set-TextValue "D8" "3.120.16"
Set-TextValue "E8" "  -0.83%  "
Set-TextValue "E9" "  -1.54%  "
Set-TextValue "E10" "  -2.83%  "
Set-TextValue "D11" "5.35"
Set-TextValue "E11" "  -0.39%  "
Set-TextValue "E12" "  -2.72%  "
Set-TextValue "E13" "  -3.21%  "
Set-TextValue "D14" "34.10"
Set-TextValue "E14" "  -4.01%  "
Set-TextValue "D15" "3.639.99"
Set-TextValue "E15" "  -0.65%  "
Set-TextValue "E16" "  +1.14%  "
Set-TextValue "D17" "63.125.06"
Set-TextValue "E17" "  -1.98%  "
Set-TextValue "D18" "3.130.97"
Set-TextValue "E18" "  -1.30%  "
Set-TextValue "D19" "6.72"
Set-TextValue "E19" "  -1.85%  "
Set-TextValue "D20" "476.96"
Set-TextValue "E20" "  +0.11%  "
Set-TextValue "D21" "14.19"
Set-TextValue "E21" "  -3.31%  "
Set-TextValue "D22" "0.699"
Set-TextValue "E22" "  -3.25%  "
Set-TextValue "D23" "7.65"
Set-TextValue "E23" "  -2.35%  "
Set-TextValue "D24" "87.57"
Set-TextValue "E24" "  +3.11%  "
Set-TextValue "D25" "13.08"
Set-TextValue "E25" "  -4.35%  "
Set-TextValue "E26" "  +0.31%  "
Set-TextValue "D27" "2.71"
Set-TextValue "E27" "  -2.95%  "
Set-TextValue "D28" "7.18"
Set-TextValue "E28" "  -3.05%  "
Set-TextValue "D29" "8.00"
Set-TextValue "E29" "  -6.96%  "
Set-TextValue "D30" "2.08"
Set-TextValue "E30" "  -0.54%  "
Set-TextValue "D31" "27.39"
Set-TextValue "E31" "  +2.85%  "
Set-TextValue "E32" "  -0.01%  "
Set-TextValue "D33" "0.109"
Set-TextValue "E33" "  -8.16%  "
Set-TextValue "E34" "  -3.89%  "
Set-TextValue "E35" "  -3.14%  "
Set-TextValue "D36" "5.86"
Set-TextValue "E36" "  -1.44%  "
Set-TextValue "D37" "52.07"
Set-TextValue "E37" "  -0.99%  "
Set-TextValue "E38" "  -3.44%  "
Set-TextValue "D39" "0.0390"
Set-TextValue "E39" "  -1.38%  "
Set-TextValue "D40" "420.90"
Set-TextValue "E40" "  -7.13%  "
Set-TextValue "E41" "  -0.80%  "
Set-TextValue "E42" "  -0.66%  "
Set-TextValue "E43" "  -10.51%  "
Set-TextValue "D44" "2.867.16"
Set-TextValue "E44" "  +0.39%  "
Set-TextValue "D45" "0.259"
Set-TextValue "E45" "  -3.00%  "
Set-TextValue "E47" "  -6.55%  "
Set-TextValue "B48" "InjectiveProtocol"
Set-TextValue "C48" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D48" "25.49"
Set-TextValue "E48" "  -3.62%  "
Set-TextValue "E49" "  -0.37%  "
Set-TextValue "B50" "ThetaToken"
Set-TextValue "C50" "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D50" "2.29"
Set-TextValue "E50" "  -5.52%  "
Set-TextValue "D51" "118.62"
Set-TextValue "E51" "  -1.23%  "
